# 4-Feb-2021, midday update.
# Adds new petty-cash entries (rows 17-26) to the "Sheet1" ledger
# (the sheet with the daily Kas Kecil / petty cash transactions,
# dimension A1:L113, tabSelected).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# --- Row 17: extra Wages Expense top-up ---
$ws.Range("D17").Formula = "=60000+300000"

# --- Row 18: A/R ---
$ws.Range("B18").Value = "A/R"
$ws.Range("C18").Formula = "=20400000+5820000+38426000"

# --- Row 19: TRANSFER BCA ---
$ws.Range("B19").Value = "TRANSFER BCA"
$ws.Range("D19").Formula = "=20400000+2574000+5000000+2262000+850000"

# --- Row 20: LPG (new item) ---
$ws.Range("B20").Value = "LPG"
$ws.Range("D20").Value = 145000

# --- Row 21: IURAN DAERAH (new item) ---
$ws.Range("B21").Value = "IURAN DAERAH"
$ws.Range("D21").Formula = "=25000"

# --- Row 22: SALES - cash/retail ---
$ws.Range("B22").Value = "SALES - cash/retail"
$ws.Range("C22").Formula = "=4392475+42166525-38426000"

# --- Row 23: SETOR KE BANK ---
$ws.Range("B23").Value = "SETOR KE BANK"
$ws.Range("D23").Formula = "=41000000"

# --- Row 24: new day, 4-Feb-2021 (serial 44231) ---
$ws.Range("A24").Value = 44231
$ws.Range("B24").Value = "Wages Expense"
$ws.Range("D24").Formula = "=60000"

# --- Row 25: TRANSFER BCA ---
$ws.Range("B25").Value = "TRANSFER BCA"
$ws.Range("D25").Formula = "=1550000+41600000"

# --- Row 26: A/R ---
$ws.Range("B26").Value = "A/R"
$ws.Range("C26").Formula = "=17240000+24360000"

# Update the view state to match: scrolled/selected around the new rows.
$ws.Range("C44").Select()
